$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 105, shifting the existing rows 105-107 down to 106-108.
$ws.Rows(105).Insert()

# Populate the newly inserted row 105 with this week's data (same market/
# category/etc. as the surrounding rows, new date + measurements).
$ws.Range("A105").Value = 6
$ws.Range("B105").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C105").Value = "Metropolitana"
$ws.Range("D105").Value = 45147
$ws.Range("D105").NumberFormat = $ws.Range("D106").NumberFormat
$ws.Range("E105").Value = 13
$ws.Range("F105").Value = 100112035
$ws.Range("G105").Value = "Bruselas (repollito)"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 410
$ws.Range("K105").Value = 17000
$ws.Range("L105").Value = 18000
$ws.Range("M105").Value = 17366
$ws.Range("N105").Value = "$/malla 15 kilos"
$ws.Range("O105").Value = "Provincia de Quillota"
$ws.Range("P105").Value = 1158
$ws.Range("Q105").Value = 15
$ws.Range("R105").Value = "Hortaliza"
